$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (customer name / card number) ---
$ws.Range("C2").Value = "Hartmut"

# Card number is a long digit string that must stay text, not be coerced
# into a Number (which would lose the leading-zero-safe formatting and
# round-trip through scientific notation) - force text format before
# writing, then restore the original cell style by re-pasting formats
# from an untouched neighbour cell that still carries it (plain
# NumberFormat changes mint a brand-new style record otherwise).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 27.10.2023"

# --- Transaction rows 6-10 (dates / description / amount) ---
$ws.Range("B6").Value = "29.10."
$ws.Range("C6").Value = "30.10."
$ws.Range("D6").Value = "KARTENZ./29.10 REWE RO"
$ws.Range("E6").Value = "115,95-"

$ws.Range("B7").Value = "31.10."
$ws.Range("C7").Value = "01.11."
$ws.Range("D7").Value = "PAYPAL FIFVBN"
$ws.Range("E7").Value = "63,34-"

$ws.Range("B8").Value = "02.11."
$ws.Range("C8").Value = "03.11."
$ws.Range("D8").Value = "KARTENZ./02.11 ALDI SUED RO"
$ws.Range("E8").Value = "139,40-"

$ws.Range("B9").Value = "04.11."
$ws.Range("C9").Value = "05.11."
$ws.Range("D9").Value = "AMAZON.DE MKTPLC EU AKZZSD"
$ws.Range("E9").Value = "85,42-"

$ws.Range("B10").Value = "08.11."
$ws.Range("C10").Value = "09.11."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "24,87-"

# --- Row 11 was blank; it now gains a new transaction line.
# Copy E10's formatting onto E11 first so it picks up the same
# right-aligned numeric style used by the other amount cells
# (row 11's E cell previously used a different, blank-row style).
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("B11").Value = "10.11."
$ws.Range("C11").Value = "11.11."
$ws.Range("D11").Value = "RECHNUNG VODAFONE GMBH 90567797"
$ws.Range("E11").Value = "38,75-"

# --- Closing balance + next billing date ---
$ws.Range("D12").Value = "KONTOSTAND AM 14.11.2023"
$ws.Range("E12").Value = "467,73-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 22.11.2023"
